# Update the "truong THPT" sample rows so the "Ma Truong" (school code) and
# "MaTpTruong" (combined province+school code) columns reflect the real
# school codes (1000 / 011000 and 1002 / 011002) instead of the old
# placeholder codes (800 / 01800 and 801 / 01800).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (STT 5500)
$ws.Range("G2").Value = "1000"     # Ma Truong
$ws.Range("B2").Value = "011000"   # MaTpTruong

# Row 3 (STT 5501)
$ws.Range("G3").Value = "1002"     # Ma Truong
$ws.Range("B3").Value = "011002"   # MaTpTruong

# Match the author's final selection in the saved workbook.
$ws.Range("C7").Select()
